$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts GFA/IMF/OECD columns right)
$ws.Columns("C").Insert()

# Header for the new column
$ws.Range("C1").Value = "M_PL"

# New profit/loss values for the new column, one per group row
$ws.Range("C2").Value = 139352925129
$ws.Range("C3").Value = 1017808846331
$ws.Range("C4").Value = 41087099249
$ws.Range("C5").Value = 47942268783
$ws.Range("C6").Value = 959424197928
$ws.Range("C7").Value = 38137333219
$ws.Range("C8").Value = 91572151625
$ws.Range("C9").Value = 1068114426119
$ws.Range("C10").Value = 200262379816
$ws.Range("C11").Value = 867852046303
